$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 4390648
$ws.Range("I9").Value = 9259732
$ws.Range("J9").Value = 8472
$ws.Range("K9").Value = 9259732
$ws.Range("L9").Value = 8472
$ws.Range("M9").Value = -9259563
$ws.Range("N9").Value = -8810
# Row 41
$ws.Range("H41").Value = 1159.8889
$ws.Range("I41").Value = 632.0769
$ws.Range("J41").Value = 1650
$ws.Range("K41").Value = 632.0769
$ws.Range("L41").Value = 1650
$ws.Range("M41").Value = -192.0769
$ws.Range("N41").Value = -2530
# Row 137
$ws.Range("H137").Value = 7144341.5
$ws.Range("I137").Value = 1799.6666
$ws.Range("J137").Value = 9092307
$ws.Range("K137").Value = 5398.9998
$ws.Range("L137").Value = 27276921
$ws.Range("M137").Value = -2848.9998
$ws.Range("N137").Value = -27282021

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 3374.25
$ws.Range("I16").Value = 3245
$ws.Range("J16").Value = 3503.5
$ws.Range("K16").Value = 3245
$ws.Range("L16").Value = 3503.5
$ws.Range("M16").Value = -2958
$ws.Range("N16").Value = -4077.5
# Row 30
$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 500
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -350
$ws.Range("N30").ClearContents()
# Row 32
$ws.Range("H32").Value = 246620.05
$ws.Range("I32").Value = 287026.62
$ws.Range("K32").Value = 287026.62
$ws.Range("M32").Value = -286739.62
# Row 35
$ws.Range("H35").Value = 3845.3333
$ws.Range("I35").Value = 3845.3333
$ws.Range("K35").Value = 3845.3333
$ws.Range("M35").Value = -3439.3333
# Row 45
$ws.Range("H45").Value = 48631.816
$ws.Range("I45").Value = 65320.125
$ws.Range("K45").Value = 65320.125
$ws.Range("M45").Value = -64943.125
# Row 61
$ws.Range("H61").Value = 786901.5
$ws.Range("I61").Value = 2400.0186
$ws.Range("K61").Value = 2400.0186
$ws.Range("M61").Value = -2188.0186
# Row 74
$ws.Range("H74").Value = 639746.9399999999
$ws.Range("I74").Value = 1835.1765
$ws.Range("K74").Value = 1835.1765
$ws.Range("M74").Value = -961.1765
# Row 77
$ws.Range("H77").Value = 639746.9399999999
$ws.Range("I77").Value = 1835.1765
$ws.Range("K77").Value = 9175.8825
$ws.Range("M77").Value = -4807.8825
# Row 132
$ws.Range("H132").Value = 2295.2295
$ws.Range("I132").Value = 2295.2295
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6885.6885
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4355.6885
$ws.Range("N132").ClearContents()
# Row 135
$ws.Range("H135").Value = 95185.42999999999
$ws.Range("J135").Value = 96883
$ws.Range("L135").Value = 96883
$ws.Range("N135").Value = -107023
# Row 136
$ws.Range("H136").Value = 786901.5
$ws.Range("I136").Value = 2400.0186
$ws.Range("K136").Value = 7200.0558
$ws.Range("M136").Value = -4650.0558

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 1007.0625
$ws.Range("I22").Value = 929.5
$ws.Range("J22").Value = 1550
$ws.Range("K22").Value = 929.5
$ws.Range("L22").Value = 1550
$ws.Range("M22").Value = -756.5
$ws.Range("N22").Value = -1896
# Row 135
$ws.Range("H135").Value = 60551.668
$ws.Range("J135").Value = 60551.668
$ws.Range("L135").Value = 60551.668
$ws.Range("N135").Value = -70691.66800000001

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 1287072.2
$ws.Range("I3").Value = 1287072.2
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1287072.2
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1286959.2
$ws.Range("N3").ClearContents()
# Row 22
$ws.Range("H22").Value = 2024.75
$ws.Range("I22").Value = 1999.5
$ws.Range("K22").Value = 1999.5
$ws.Range("M22").Value = -1649.5
# Row 31
$ws.Range("H31").Value = 3130.8
$ws.Range("I31").Value = 2803.2083
$ws.Range("J31").Value = 3622.1875
$ws.Range("K31").Value = 2803.2083
$ws.Range("L31").Value = 3622.1875
$ws.Range("M31").Value = -2508.2083
$ws.Range("N31").Value = -4212.1875
# Row 34
$ws.Range("H34").Value = 3130.8
$ws.Range("I34").Value = 2803.2083
$ws.Range("J34").Value = 3622.1875
$ws.Range("K34").Value = 2803.2083
$ws.Range("L34").Value = 3622.1875
$ws.Range("M34").Value = -2601.2083
$ws.Range("N34").Value = -4026.1875
# Row 62
$ws.Range("H62").Value = 12001.667
$ws.Range("I62").Value = 13002.5
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 13002.5
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -12378.5
$ws.Range("N62").Value = -11248
# Row 65
$ws.Range("H65").Value = 12001.667
$ws.Range("I65").Value = 13002.5
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 65012.5
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -61892.5
$ws.Range("N65").Value = -56240
# Row 107
$ws.Range("H107").Value = 2576.8823
$ws.Range("I107").Value = 2182.4546
$ws.Range("J107").Value = 3300
$ws.Range("K107").Value = 2182.4546
$ws.Range("L107").Value = 3300
$ws.Range("M107").Value = -262.4546
$ws.Range("N107").Value = -7140
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 134
$ws.Range("H134").Value = 1373.64
$ws.Range("I134").Value = 1239.7646
$ws.Range("K134").Value = 3719.2938
$ws.Range("M134").Value = -1184.2938

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1638.1364
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1638.1364
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 4914.4092
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -14994.4092

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Range("H41").Value = 7149.1
$ws.Range("I41").Value = 7149.1
$ws.Range("K41").Value = 7149.1
$ws.Range("M41").Value = -6794.1
# Row 70
$ws.Range("H70").Value = 9220.166999999999
$ws.Range("I70").Value = 9000
$ws.Range("J70").Value = 9440.333000000001
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 9440.333000000001
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -9980.333000000001
# Row 73
$ws.Range("H73").Value = 9220.166999999999
$ws.Range("I73").Value = 9000
$ws.Range("J73").Value = 9440.333000000001
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 9440.333000000001
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -11312.333
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
# Row 117
$ws.Range("H117").Value = 49750
$ws.Range("J117").Value = 49750
$ws.Range("L117").Value = 49750
$ws.Range("N117").Value = -56634
# Row 122
$ws.Range("H122").Value = 1978.2368
$ws.Range("I122").Value = 1723.6666
$ws.Range("K122").Value = 5170.9998
$ws.Range("M122").Value = -2720.9998
# Row 132
$ws.Range("H132").Value = 992866.6
$ws.Range("I132").Value = 11226.083
$ws.Range("J132").Value = 1834272.8
$ws.Range("K132").Value = 33678.249
$ws.Range("L132").Value = 5502818.4
$ws.Range("M132").Value = -31148.249
$ws.Range("N132").Value = -5507878.4
# Row 140
$ws.Range("H140").Value = 102547.6
$ws.Range("J140").Value = 105085.22
$ws.Range("L140").Value = 105085.22
$ws.Range("N140").Value = -115445.22

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 122300000
$ws.Range("J20").Value = 366700000
$ws.Range("L20").Value = 366700000
$ws.Range("N20").Value = -366700452
# Row 22
$ws.Range("H22").Value = 3865.6667
$ws.Range("I22").Value = 1271.8182
$ws.Range("K22").Value = 1271.8182
$ws.Range("M22").Value = -976.8181999999999
# Row 27
$ws.Range("H27").Value = 3865.6667
$ws.Range("I27").Value = 1271.8182
$ws.Range("K27").Value = 1271.8182
$ws.Range("M27").Value = -1164.8182
# Row 42
$ws.Range("H42").Value = 96633
$ws.Range("I42").Value = 100000
$ws.Range("J42").Value = 89899
$ws.Range("K42").Value = 100000
$ws.Range("L42").Value = 89899
$ws.Range("M42").Value = -99437
$ws.Range("N42").Value = -91025
# Row 43
$ws.Range("H43").Value = 4281389
$ws.Range("J43").Value = 6828636.5
$ws.Range("L43").Value = 6828636.5
$ws.Range("N43").Value = -6829022.5
# Row 49
$ws.Range("H49").Value = 96633
$ws.Range("I49").Value = 100000
$ws.Range("J49").Value = 89899
$ws.Range("K49").Value = 100000
$ws.Range("L49").Value = 89899
$ws.Range("M49").Value = -99853
$ws.Range("N49").Value = -90193
# Row 99
$ws.Range("H99").Value = 50324.832
$ws.Range("I99").Value = 36718
$ws.Range("K99").Value = 36718
$ws.Range("M99").Value = -33723
# Row 132
$ws.Range("H132").Value = 3082.697
$ws.Range("I132").Value = 2797.3914
$ws.Range("K132").Value = 8392.174199999999
$ws.Range("M132").Value = -5862.174199999999
# Row 136
$ws.Range("H136").Value = 2137.4773
$ws.Range("I136").Value = 2597.68
$ws.Range("J136").Value = 1954.8572
$ws.Range("K136").Value = 7793.039999999999
$ws.Range("L136").Value = 5864.571599999999
$ws.Range("M136").Value = -5243.039999999999
$ws.Range("N136").Value = -10964.5716

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 2916.4546
$ws.Range("I62").Value = 2864.2222
$ws.Range("K62").Value = 2864.2222
$ws.Range("M62").Value = -2240.2222
# Row 65
$ws.Range("H65").Value = 2916.4546
$ws.Range("I65").Value = 2864.2222
$ws.Range("K65").Value = 14321.111
$ws.Range("M65").Value = -11201.111
# Row 122
$ws.Range("H122").Value = 1750.6364
$ws.Range("J122").Value = 4650
$ws.Range("L122").Value = 13950
$ws.Range("N122").Value = -18850
# Row 132
$ws.Range("H132").Value = 2403.6047
$ws.Range("I132").Value = 1952.8889
$ws.Range("J132").Value = 3164.1875
$ws.Range("K132").Value = 5858.6667
$ws.Range("L132").Value = 9492.5625
$ws.Range("M132").Value = -3328.6667
$ws.Range("N132").Value = -14552.5625
